# Insert a new timestamp column before the "nom" / "url_produit" columns.
# Before: ... BY (last price snapshot), BZ = "nom", CA = "url_produit"
# After:  ... BY (last price snapshot), BZ = new timestamp snapshot,
#             CA = "nom", CB = "url_produit"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift BZ:CA (and everything to their right) one column to the right,
# inserting a brand-new blank column at BZ. xlShiftToRight = -4161.
$ws.Range("BZ1:BZ206").Insert(-4161)

# New header cell for the inserted timestamp column.
$ws.Range("BZ1").Value = "2026-01-31 06:24:29"

# Rows 2-80 already had a numeric price snapshot recorded in BY (the most
# recent prior timestamp column); duplicate that same price into the new
# BZ column for this snapshot.
for ($r = 2; $r -le 80; $r++) {
    $priceValue = $ws.Cells.Item($r, 77).Value()
    $ws.Cells.Item($r, 78).Value = $priceValue
}

# Rows 81-206 had no price recorded for this run (BY is blank), so the new
# BZ cell for those rows is left blank as well - which is already the case
# right after the column insert, so nothing else to do for them.
